$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header tweaks
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 data tweaks
$ws.Range("B2").Value = 16.626316213124994
$ws.Range("C2").Value = 20.806062633125009
$ws.Range("D2").Value = 24.913831213124979
$ws.Range("E2").Value = 24.048378633124969

# Row 3 data tweaks
$ws.Range("B3").Value = 18.431477713124991
$ws.Range("C3").Value = 26.924859802499952
$ws.Range("D3").Value = 21.178336213124965
$ws.Range("E3").Value = 23.074525552499949

# Selection now covers only the edited block
$ws.Range("B1:E3").Select()
